# Update loading_percent results for case with 380 kV (rows 2-25)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("B2").Value = 19.26372206575458
$ws.Range("D2").Value = 9.44126855360463
$ws.Range("E2").Value = 14.92188802172695
$ws.Range("F2").Value = 30.90481186196838
$ws.Range("G2").Value = 3.63964322550747
$ws.Range("I2").Value = 24.39325557877157
$ws.Range("J2").Value = 10.7612574448854
$ws.Range("L2").Value = 13.25964699502668
$ws.Range("O2").Value = 23.54983101980319

# Row 3
$ws.Range("B3").Value = 18.65906230001512
$ws.Range("D3").Value = 9.44801843480195
$ws.Range("E3").Value = 14.84745460409717
$ws.Range("F3").Value = 31.05686485496661
$ws.Range("G3").Value = 3.642218384008967
$ws.Range("I3").Value = 24.57077162536456
$ws.Range("J3").Value = 10.72634691654375
$ws.Range("L3").Value = 12.90651743915028
$ws.Range("O3").Value = 23.62124284357227

# Row 4
$ws.Range("B4").Value = 18.27727883939972
$ws.Range("D4").Value = 9.453042502047909
$ws.Range("E4").Value = 14.80311576290471
$ws.Range("F4").Value = 31.16070733145129
$ws.Range("G4").Value = 3.643883505373849
$ws.Range("I4").Value = 24.6869950528856
$ws.Range("J4").Value = 10.70584606458029
$ws.Range("L4").Value = 12.6845587708625
$ws.Range("O4").Value = 23.6727507468271

# Row 5
$ws.Range("B5").Value = 18.1192482056706
$ws.Range("D5").Value = 9.455311927697485
$ws.Range("E5").Value = 14.78540343950279
$ws.Range("F5").Value = 31.20564474879775
$ws.Range("G5").Value = 3.644583239287857
$ws.Range("I5").Value = 24.73617263952232
$ws.Range("J5").Value = 10.69773212261343
$ws.Range("L5").Value = 12.59292680557041
$ws.Range("O5").Value = 23.69565709603989

# Row 6
$ws.Range("B6").Value = 18.09286533609602
$ws.Range("D6").Value = 9.455702202771713
$ws.Range("E6").Value = 14.78248420949609
$ws.Range("F6").Value = 31.21326442091591
$ws.Range("G6").Value = 3.644700711059431
$ws.Range("I6").Value = 24.74444814591033
$ws.Range("J6").Value = 10.69639947714347
$ws.Range("L6").Value = 12.5776433373643
$ws.Range("O6").Value = 23.69957612613725

# Row 7
$ws.Range("B7").Value = 18.27515723665171
$ws.Range("D7").Value = 9.453072207962837
$ws.Range("E7").Value = 14.8028754292277
$ws.Range("F7").Value = 31.16130278169872
$ws.Range("G7").Value = 3.643892856366768
$ws.Range("I7").Value = 24.68765093015605
$ws.Range("J7").Value = 10.70573565710859
$ws.Range("L7").Value = 12.68332762612074
$ws.Range("O7").Value = 23.67305192377548

# Row 8
$ws.Range("B8").Value = 19.05753027055244
$ws.Range("D8").Value = 9.443413888211984
$ws.Range("E8").Value = 14.89594841149201
$ws.Range("F8").Value = 30.95505549443504
$ws.Range("G8").Value = 3.64051375577897
$ws.Range("I8").Value = 24.45296177578345
$ws.Range("J8").Value = 10.74902963877765
$ws.Range("L8").Value = 13.13901719604758
$ws.Range("O8").Value = 23.57285805813816

# Row 9
$ws.Range("B9").Value = 20.50060006712956
$ws.Range("D9").Value = 9.431413933261917
$ws.Range("E9").Value = 15.08868842623669
$ws.Range("F9").Value = 30.63442004185876
$ws.Range("G9").Value = 3.634550365250472
$ws.Range("I9").Value = 24.05017866258719
$ws.Range("J9").Value = 10.84110105752239
$ws.Range("L9").Value = 13.98744571750455
$ws.Range("O9").Value = 23.43757743274715

# Row 10
$ws.Range("B10").Value = 21.49611516416753
$ws.Range("D10").Value = 9.426775205131078
$ws.Range("E10").Value = 15.23571218448944
$ws.Range("F10").Value = 30.4508125409377
$ws.Range("G10").Value = 3.630568777750028
$ws.Range("I10").Value = 23.78940782248305
$ws.Range("J10").Value = 10.91278660887092
$ws.Range("L10").Value = 14.577852703279
$ws.Range("O10").Value = 23.37601930644954

# Row 11
$ws.Range("B11").Value = 21.93333429452536
$ws.Range("D11").Value = 9.425560386848726
$ws.Range("E11").Value = 15.30359478487996
$ws.Range("F11").Value = 30.37874599455934
$ws.Range("G11").Value = 3.628843291424489
$ws.Range("I11").Value = 23.6784438219205
$ws.Range("J11").Value = 10.94620211510115
$ws.Range("L11").Value = 14.83830426521829
$ws.Range("O11").Value = 23.3563266658234

# Row 12
$ws.Range("B12").Value = 22.09653267906503
$ws.Range("D12").Value = 9.425228095790072
$ws.Range("E12").Value = 15.32942866765696
$ws.Range("F12").Value = 30.35311780918817
$ws.Range("G12").Value = 3.628202153918678
$ws.Range("I12").Value = 23.63752985457988
$ws.Range("J12").Value = 10.95896506193836
$ws.Range("L12").Value = 14.93568913741624
$ws.Range("O12").Value = 23.35007124102467

# Row 13
$ws.Range("B13").Value = 22.06149198373073
$ws.Range("D13").Value = 9.425293996039803
$ws.Range("E13").Value = 15.32385942748462
$ws.Range("F13").Value = 30.35856316421327
$ws.Range("G13").Value = 3.628339689819228
$ws.Range("I13").Value = 23.64629216853503
$ws.Range("J13").Value = 10.95621158154032
$ws.Range("L13").Value = 14.91477192193021
$ws.Range("O13").Value = 23.35136491282411

# Row 14
$ws.Range("B14").Value = 21.94680878571823
$ws.Range("D14").Value = 9.42553049498572
$ws.Range("E14").Value = 15.3057176746489
$ws.Range("F14").Value = 30.376604162363
$ws.Range("G14").Value = 3.628790299184228
$ws.Range("I14").Value = 23.67505562218347
$ws.Range("J14").Value = 10.94724997964661
$ws.Range("L14").Value = 14.84634147583252
$ws.Range("O14").Value = 23.35578790801531

# Row 15
$ws.Range("B15").Value = 21.87625051300423
$ws.Range("D15").Value = 9.425691961611831
$ws.Range("E15").Value = 15.29462155680861
$ws.Range("F15").Value = 30.38787163195045
$ws.Range("G15").Value = 3.629067905994257
$ws.Range("I15").Value = 23.69281818220096
$ws.Range("J15").Value = 10.94177476580511
$ws.Range("L15").Value = 14.80426197560788
$ws.Range("O15").Value = 23.35865380818326

# Row 16
$ws.Range("B16").Value = 21.46721670720195
$ws.Range("D16").Value = 9.426872526841997
$ws.Range("E16").Value = 15.23129469621326
$ws.Range("F16").Value = 30.45575404502528
$ws.Range("G16").Value = 3.630683262426405
$ws.Range("I16").Value = 23.79681408755953
$ws.Range("J16").Value = 10.91061847230105
$ws.Range("L16").Value = 14.5606615676847
$ws.Range("O16").Value = 23.37747408847546

# Row 17
$ws.Range("B17").Value = 21.21219275741212
$ws.Range("D17").Value = 9.427825371956663
$ws.Range("E17").Value = 15.19269091385699
$ws.Range("F17").Value = 30.50034255587807
$ws.Range("G17").Value = 3.631696149413518
$ws.Range("I17").Value = 23.86257732045596
$ws.Range("J17").Value = 10.89170671912705
$ws.Range("L17").Value = 14.40908421080306
$ws.Range("O17").Value = 23.39115328180627

# Row 18
$ws.Range("B18").Value = 21.0640431868133
$ws.Range("D18").Value = 9.428457761278773
$ws.Range("E18").Value = 15.17058218375061
$ws.Range("F18").Value = 30.52706629156053
$ws.Range("G18").Value = 3.632286811014186
$ws.Range("I18").Value = 23.90112352278157
$ws.Range("J18").Value = 10.88090538798995
$ws.Range("L18").Value = 14.32113986897913
$ws.Range("O18").Value = 23.39980306369427

# Row 19
$ws.Range("B19").Value = 21.01363420234014
$ws.Range("D19").Value = 9.428686393221742
$ws.Range("E19").Value = 15.16311336473848
$ws.Range("F19").Value = 30.53629912511126
$ws.Range("G19").Value = 3.632488187902081
$ws.Range("I19").Value = 23.91429831013372
$ws.Range("J19").Value = 10.87726153093153
$ws.Range("L19").Value = 14.29123507781966
$ws.Range("O19").Value = 23.40286579968599

# Row 20
$ws.Range("B20").Value = 21.23949316535925
$ws.Range("D20").Value = 9.427715218607343
$ws.Range("E20").Value = 15.19679061855543
$ws.Range("F20").Value = 30.49548439966149
$ws.Range("G20").Value = 3.631587490522973
$ws.Range("I20").Value = 23.85550207093234
$ws.Range("J20").Value = 10.89371206759197
$ws.Range("L20").Value = 14.4252992022814
$ws.Range("O20").Value = 23.38961613608051

# Row 21
$ws.Range("B21").Value = 21.98055912819876
$ws.Range("D21").Value = 9.425457571061036
$ws.Range("E21").Value = 15.31104299314706
$ws.Range("F21").Value = 30.37125987042855
$ws.Range("G21").Value = 3.628657611892865
$ws.Range("I21").Value = 23.66657705801232
$ws.Range("J21").Value = 10.94987930824927
$ws.Range("L21").Value = 14.86647543003813
$ws.Range("O21").Value = 23.3544561027344

# Row 22
$ws.Range("B22").Value = 22.45104674045034
$ws.Range("D22").Value = 9.424726113274868
$ws.Range("E22").Value = 15.38645350758314
$ws.Range("F22").Value = 30.29976475318678
$ws.Range("G22").Value = 3.62681423672851
$ws.Range("I22").Value = 23.54955045544525
$ws.Range("J22").Value = 10.98722133153234
$ws.Range("L22").Value = 15.14753815725222
$ws.Range("O22").Value = 23.33848420472087

# Row 23
$ws.Range("B23").Value = 22.20124046217121
$ws.Range("D23").Value = 9.425048767746471
$ws.Range("E23").Value = 15.34614300120512
$ws.Range("F23").Value = 30.33703156834627
$ws.Range("G23").Value = 3.627791562660112
$ws.Range("I23").Value = 23.61141851248851
$ws.Range("J23").Value = 10.96723543283572
$ws.Range("L23").Value = 14.99821758653654
$ws.Range("O23").Value = 23.3463655546051

# Row 24
$ws.Range("B24").Value = 21.22715541652073
$ws.Range("D24").Value = 9.42776475539068
$ws.Range("E24").Value = 15.19493687578926
$ws.Range("F24").Value = 30.49767738027778
$ws.Range("G24").Value = 3.631636589224754
$ws.Range("I24").Value = 23.85869849250696
$ws.Range("J24").Value = 10.89280522698785
$ws.Range("L24").Value = 14.41797089236453
$ws.Range("O24").Value = 23.3903086335442

# Row 25
$ws.Range("B25").Value = 20.1209763108901
$ws.Range("D25").Value = 9.433922674081154
$ws.Range("E25").Value = 15.03553563617304
$ws.Range("F25").Value = 30.71209642148068
$ws.Range("G25").Value = 3.636093104365568
$ws.Range("I25").Value = 24.15297950586439
$ws.Range("J25").Value = 10.81545910600688
$ws.Range("L25").Value = 13.76334623740298
$ws.Range("O25").Value = 23.46756566391307
